# ADD: graph variations for base models
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("KG Default (Balanced 0.5-0.5)", 0.4728205128205129, 0.7076923076923077, 3.811671018600464),
    @("KG Artem V1 Full (Balanced 0.5-0.5)", 0.3069230769230769, 0.6307692307692307, 0.7538149356842041),
    @("KG Artem V2 Purpose+Insights (Balanced 0.5-0.5)", 0.3430769230769231, 0.6615384615384615, 0.740900993347168),
    @("KG Artem V3 Full (Balanced 0.5-0.5)", 0.3069230769230769, 0.6307692307692307, 0.7330317497253418),
    @("KG Artem V4 Minimal (Balanced 0.5-0.5)", 0.3374358974358974, 0.6307692307692307, 0.7192442417144775),
    @("KG Default (BM25 0.7-0.3)", 0.4123076923076923, 0.6923076923076923, 0.7164630889892578),
    @("KG Default (Vector 0.3-0.7)", 0.4587179487179487, 0.7230769230769231, 0.7145400047302246)
)

$startRow = 29
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
